$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text (Volume NN -> NN, issue number)
$ws.Range("A8").Characters(21, 2).Text = "30"

# Update report date range text (week covering dates)
$ws.Range("C9").Characters(27, 9).Text = "7/24/2023"
$ws.Range("C9").Characters(47, 9).Text = "7/30/2023"

# Update weekly crime statistics table (rows 14-30)
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = -40
$ws.Range("I14").Value = 35
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 16.666666666666
$ws.Range("L14").Value = -2.777777777777
$ws.Range("M14").Value = 16.666666666666
$ws.Range("N14").Value = -82.233502538071
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 200
$ws.Range("F15").Value = 11
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = -31.25
$ws.Range("I15").Value = 77
$ws.Range("J15").Value = 111
$ws.Range("K15").Value = -30.630630630630
$ws.Range("L15").Value = -2.531645569620
$ws.Range("M15").Value = -33.620689655172
$ws.Range("N15").Value = -72.695035460992
$ws.Range("C16").Value = 44
$ws.Range("D16").Value = 49
$ws.Range("E16").Value = -10.204081632653
$ws.Range("F16").Value = 163
$ws.Range("H16").Value = -9.944751381215
$ws.Range("I16").Value = 1081
$ws.Range("J16").Value = 1191
$ws.Range("K16").Value = -9.235936188077
$ws.Range("L16").Value = 12.486992715920
$ws.Range("M16").Value = -18.353474320241
$ws.Range("N16").Value = -80.223197950969
$ws.Range("C17").Value = 72
$ws.Range("D17").Value = 79
$ws.Range("E17").Value = -8.860759493670
$ws.Range("F17").Value = 268
$ws.Range("G17").Value = 275
$ws.Range("H17").Value = -2.545454545454
$ws.Range("I17").Value = 1739
$ws.Range("J17").Value = 1806
$ws.Range("K17").Value = -3.709856035437
$ws.Range("L17").Value = 10.063291139240
$ws.Range("M17").Value = 54.991087344028
$ws.Range("N17").Value = -49.696268440844
$ws.Range("C18").Value = 29
$ws.Range("D18").Value = 28
$ws.Range("E18").Value = 3.571428571428
$ws.Range("F18").Value = 116
$ws.Range("G18").Value = 124
$ws.Range("H18").Value = -6.451612903225
$ws.Range("I18").Value = 887
$ws.Range("J18").Value = 1049
$ws.Range("K18").Value = -15.443279313632
$ws.Range("L18").Value = 12.420785804816
$ws.Range("M18").Value = 15.194805194805
$ws.Range("N18").Value = -85.871296591271
$ws.Range("C19").Value = 135
$ws.Range("D19").Value = 155
$ws.Range("E19").Value = -12.903225806451
$ws.Range("F19").Value = 513
$ws.Range("G19").Value = 564
$ws.Range("H19").Value = -9.042553191489
$ws.Range("I19").Value = 3530
$ws.Range("J19").Value = 3725
$ws.Range("K19").Value = -5.234899328859
$ws.Range("L19").Value = 28.084179970972
$ws.Range("M19").Value = 33.661491859144
$ws.Range("N19").Value = -43.798758159528
$ws.Range("C20").Value = 25
$ws.Range("D20").Value = 34
$ws.Range("E20").Value = -26.470588235294
$ws.Range("F20").Value = 115
$ws.Range("G20").Value = 122
$ws.Range("H20").Value = -5.737704918032
$ws.Range("I20").Value = 793
$ws.Range("J20").Value = 756
$ws.Range("K20").Value = 4.894179894179
$ws.Range("L20").Value = 54.580896686159
$ws.Range("M20").Value = 146.273291925466
$ws.Range("N20").Value = -85.398637451666
$ws.Range("C21").Value = 309
$ws.Range("D21").Value = 347
$ws.Range("E21").Value = -10.951008645533
$ws.Range("F21").Value = 1192
$ws.Range("G21").Value = 1292
$ws.Range("H21").Value = -7.739938080495
$ws.Range("I21").Value = 8142
$ws.Range("J21").Value = 8668
$ws.Range("K21").Value = -6.068297185048
$ws.Range("L21").Value = 21.268990169794
$ws.Range("M21").Value = 28.727272727272
$ws.Range("N21").Value = -70.275992990654
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = -60
$ws.Range("F22").Value = 15
$ws.Range("G22").Value = 29
$ws.Range("H22").Value = -48.275862068965
$ws.Range("I22").Value = 163
$ws.Range("J22").Value = 187
$ws.Range("K22").Value = -12.834224598930
$ws.Range("L22").Value = 28.346456692913
$ws.Range("M22").Value = 22.556390977443
$ws.Range("D23").Value = 35
$ws.Range("E23").Value = -20
$ws.Range("F23").Value = 108
$ws.Range("G23").Value = 113
$ws.Range("H23").Value = -4.424778761061
$ws.Range("I23").Value = 715
$ws.Range("J23").Value = 753
$ws.Range("K23").Value = -5.046480743691
$ws.Range("L23").Value = 2.729885057471
$ws.Range("M23").Value = 51.162790697674
$ws.Range("C24").Value = 292
$ws.Range("D24").Value = 321
$ws.Range("E24").Value = -9.034267912772
$ws.Range("F24").Value = 1138
$ws.Range("G24").Value = 1262
$ws.Range("H24").Value = -9.825673534072
$ws.Range("I24").Value = 8164
$ws.Range("J24").Value = 9116
$ws.Range("K24").Value = -10.443176831943
$ws.Range("L24").Value = 18.956724464519
$ws.Range("M24").Value = 57.090629209159
$ws.Range("C25").Value = 88
$ws.Range("D25").Value = 72
$ws.Range("E25").Value = 22.222222222222
$ws.Range("F25").Value = 349
$ws.Range("G25").Value = 313
$ws.Range("H25").Value = 11.501597444089
$ws.Range("I25").Value = 2659
$ws.Range("J25").Value = 2656
$ws.Range("K25").Value = 0.112951807228
$ws.Range("L25").Value = 14.958927799394
$ws.Range("M25").Value = -15.960809102402
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = -37.037037037037
$ws.Range("I26").Value = 145
$ws.Range("J26").Value = 177
$ws.Range("K26").Value = -18.079096045197
$ws.Range("L26").Value = 2.836879432624
$ws.Range("C27").Value = 14
$ws.Range("E27").Value = 27.272727272727
$ws.Range("F27").Value = 58
$ws.Range("G27").Value = 47
$ws.Range("H27").Value = 23.404255319148
$ws.Range("I27").Value = 359
$ws.Range("J27").Value = 401
$ws.Range("K27").Value = -10.473815461346
$ws.Range("L27").Value = -1.643835616438
$ws.Range("C28").Value = 11
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 175
$ws.Range("F28").Value = 20
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 96
$ws.Range("J28").Value = 119
$ws.Range("K28").Value = -19.327731092437
$ws.Range("L28").Value = -34.246575342465
$ws.Range("M28").Value = -15.789473684210
$ws.Range("N28").Value = -80
$ws.Range("C29").Value = 9
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 350
$ws.Range("F29").Value = 18
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = 63.636363636363
$ws.Range("I29").Value = 86
$ws.Range("J29").Value = 98
$ws.Range("K29").Value = -12.244897959183
$ws.Range("L29").Value = -34.351145038167
$ws.Range("M29").Value = -13.131313131313
$ws.Range("N29").Value = -80.454545454545
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 166.666666666667
$ws.Range("I30").Value = 42
$ws.Range("J30").Value = 63
$ws.Range("K30").Value = -33.333333333333
$ws.Range("L30").Value = -12.5
